# Generate Report for handback
# Updates the zh-cn and de-de localization-status worksheets: marks the two
# source files as handed back (in sync with en-US), populates the "Latest
# Target File" / "Latest Handback File" columns (with matching hyperlinks),
# and stamps the handback datetime.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# Cornflower-blue hyperlink font color (matches the existing HyperLink style,
# BGR-encoded the way the Font.Color OLE property expects it).
$hyperlinkColor = 15570276

function Set-HandbackRow($ws, $rowNum, $sourceDisplay, $sourceUrl, $handoffDisplay, $handoffUrl, $handbackDateTime) {

    # Status -> "Handed back: in sync with en-US"
    $ws.Range("B$rowNum").Value = $statusHandedBack

    # Latest Target File (E) - same file/link as the source markdown (A)
    $eCell = $ws.Range("E$rowNum")
    $ws.Hyperlinks.Add($eCell, $sourceUrl, [Type]::Missing, [Type]::Missing, $sourceDisplay)
    $eCell.Font.Underline = 2
    $eCell.Font.Color = $hyperlinkColor

    # Latest Handback File (F) - same file/link as the handoff xlf (C)
    $fCell = $ws.Range("F$rowNum")
    $ws.Hyperlinks.Add($fCell, $handoffUrl, [Type]::Missing, [Type]::Missing, $handoffDisplay)
    $fCell.Font.Underline = 2
    $fCell.Font.Color = $hyperlinkColor

    # Latest Handback DateTime (G)
    $ws.Range("G$rowNum").Value = $handbackDateTime
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-HandbackRow $wsZh "2" `
    "6378e7a4-4c34-4a45-987c-e3baeb12303f.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/f07c4622a57e91c9827cc8d5c2688faf86f6f507/e2e/6378e7a4-4c34-4a45-987c-e3baeb12303f.md" `
    "6378e7a4-4c34-4a45-987c-e3baeb12303f.721f826828e83cc1026f124e830d1456e79c0502.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/159749a5cd8653792a07c0b7e8e510615316b7ff/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/6378e7a4-4c34-4a45-987c-e3baeb12303f.721f826828e83cc1026f124e830d1456e79c0502.zh-cn.xlf" `
    "2016-01-11 13:09:44"

Set-HandbackRow $wsZh "3" `
    "a5482964-c5c3-4d41-b138-4b1767023724.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/f07c4622a57e91c9827cc8d5c2688faf86f6f507/e2e/a5482964-c5c3-4d41-b138-4b1767023724.md" `
    "a5482964-c5c3-4d41-b138-4b1767023724.ed7dc443603a6858142309d72ed099bccfe5843f.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/159749a5cd8653792a07c0b7e8e510615316b7ff/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a5482964-c5c3-4d41-b138-4b1767023724.ed7dc443603a6858142309d72ed099bccfe5843f.zh-cn.xlf" `
    "2016-01-11 13:09:44"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

Set-HandbackRow $wsDe "2" `
    "6378e7a4-4c34-4a45-987c-e3baeb12303f.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/f07c4622a57e91c9827cc8d5c2688faf86f6f507/e2e/6378e7a4-4c34-4a45-987c-e3baeb12303f.md" `
    "6378e7a4-4c34-4a45-987c-e3baeb12303f.721f826828e83cc1026f124e830d1456e79c0502.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea141e3a620bc6ff24ef1326969f9cf6d93508a5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/6378e7a4-4c34-4a45-987c-e3baeb12303f.721f826828e83cc1026f124e830d1456e79c0502.de-de.xlf" `
    "2016-01-11 13:10:16"

Set-HandbackRow $wsDe "3" `
    "a5482964-c5c3-4d41-b138-4b1767023724.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/f07c4622a57e91c9827cc8d5c2688faf86f6f507/e2e/a5482964-c5c3-4d41-b138-4b1767023724.md" `
    "a5482964-c5c3-4d41-b138-4b1767023724.ed7dc443603a6858142309d72ed099bccfe5843f.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea141e3a620bc6ff24ef1326969f9cf6d93508a5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a5482964-c5c3-4d41-b138-4b1767023724.ed7dc443603a6858142309d72ed099bccfe5843f.de-de.xlf" `
    "2016-01-11 13:10:16"

Write-Output "handback report generated"
